$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("OM_Exp1")
$ws2 = $wb.Worksheets.Item("OM_Exp2")

# ---- Rename the Growth_M labels on OM_Exp1 (10 -> 15) ----
$ws1.Range("A3").Value = "Growth_M (15,0)"
$ws1.Range("A4").Value = "Growth_M (0,15)"
$ws1.Range("A5").Value = "Growth_M (15,15)"
$ws1.Range("A7").Value = "Growth_M (15,30)"
$ws1.Range("A8").Value = "Growth_M (30,15)"

# ---- Update the vonB growth/natmort control factors: 0.9 -> 0.85 ----
# OM_Exp1
$ws1.Range("B3").Value = 0.85
$ws1.Range("C4").Value = 0.85
$ws1.Range("B5").Value = 0.85
$ws1.Range("C5").Value = 0.85
$ws1.Range("B7").Value = 0.85
$ws1.Range("C8").Value = 0.85

# OM_Exp2
$ws2.Range("C2").Value = 0.85
$ws2.Range("D2").Value = 0.85
$ws2.Range("C3").Value = 0.85
$ws2.Range("D3").Value = 0.85
$ws2.Range("C4").Value = 0.85
$ws2.Range("D4").Value = 0.85
$ws2.Range("C5").Value = 0.85
$ws2.Range("D5").Value = 0.85
$ws2.Range("D6").Value = 0.85
$ws2.Range("C7").Value = 0.85

# ---- Update sheet selections / active tab ----
# OM_Exp1 becomes inactive, selection moves to B3
$ws1.Range("B3").Select()

# OM_Exp2 becomes the active tab, selection moves to C7
$ws2.Activate()
$ws2.Range("C7").Select()
